# Update Andrew Case's time card by +30 minutes (84h 45m -> 85h 15m)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "85h 15m"
$ws.Range("B4").Select()
